$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These result cells hold plain numeric-looking text ("7.10" -> "7.1", etc.)
# in the original sheet. Force text format first so Excel keeps assigning a
# string (preserving values like trailing ".0") instead of silently
# converting the cell to a number.
$numericResultCells = @("C2","C3","C4","C5","C6","C7","C8","C9","C10","C11","C12","C15","C16","C17","C21","C22")
foreach ($addr in $numericResultCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Column C: numeric results - strip trailing zero from two-decimal values
$ws.Range("C2").Value = "7.1"
$ws.Range("C3").Value = "1.6"
$ws.Range("C4").Value = "4.5"
$ws.Range("C5").Value = "1.0"
$ws.Range("C6").Value = "0.0"
$ws.Range("C7").Value = "0.0"
$ws.Range("C8").Value = "23.2"
$ws.Range("C9").Value = "61.9"
$ws.Range("C10").Value = "14.6"
$ws.Range("C11").Value = "0.2"
$ws.Range("C12").Value = "0.1"
$ws.Range("C15").Value = "38.0"
$ws.Range("C16").Value = "76.5"
$ws.Range("C17").Value = "25.6"
$ws.Range("C21").Value = "6.5"
$ws.Range("C22").Value = "17.9"

# Column B: rename test item labels
$ws.Range("B4").Value = "淋巴细胞比率"
$ws.Range("B5").Value = "单核细胞数"
$ws.Range("B6").Value = "嗜酸性粒细胞数"
$ws.Range("B22").Value = "血小板分布宽度"

# Column D: reference ranges - remove stray leading markers / trailing units and normalize dash
$ws.Range("D4").Value = "0.40-4.40"
$ws.Range("D5").Value = "0.00-0.80"
$ws.Range("D8").Value = "50.00-70.00"
$ws.Range("D9").Value = "20.00-40.00"
$ws.Range("D10").Value = "0.00-9.00"
$ws.Range("D14").Value = "110-170"
$ws.Range("D17").Value = "127.00-36.00"
$ws.Range("D18").Value = "320-360"
$ws.Range("D21").Value = "5.00-10.00"
$ws.Range("D22").Value = "9.00-20.00"
